$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aggiornamento 15, 16, 17 marzo: nuove righe giornaliere (date seriali 44301-44303)
$data = @(
    @(44301, 2, 25, 253.4982762117218),
    @(44302, 2, 19, 192.6586899209085),
    @(44303, 3, 16, 162.2388967755019)
)

$lastRow = 226
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $lastRow + 1 + $i
    $row = $data[$i]

    # Copy the formatting from the row above (date style with border/bold/center/top)
    $ws.Range("A$($r - 1)").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
